$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.865.72"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.415.84"
$ws.Range("E3").Value = "  +0.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.68"
$ws.Range("E5").Value = "  +0.77%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.32"
$ws.Range("E6").Value = "  -0.50%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.43%  "

# Row 9
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("E10").Value = "  -1.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.19"
$ws.Range("E11").Value = "  -4.09%  "

# Row 12
$ws.Range("E12").Value = "  -0.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.81"
$ws.Range("E13").Value = "  -1.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000171"
$ws.Range("E14").Value = "  -1.62%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.851.80"
$ws.Range("E15").Value = "  +0.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.839.78"
$ws.Range("E16").Value = "  -0.33%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.421.44"
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
$ws.Range("E18").Value = "  +1.18%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "322.57"
$ws.Range("E19").Value = "  -0.45%  "

# Row 20
$ws.Range("E20").Value = "  +0.77%  "

# Row 21
$ws.Range("E21").Value = "  -1.88%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.43"
$ws.Range("E23").Value = "  +1.65%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.71"
$ws.Range("E24").Value = "  -0.72%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.67"
$ws.Range("E25").Value = "  -5.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "551.96"
$ws.Range("E26").Value = "  -5.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.537.51"
$ws.Range("E27").Value = "  +0.31%  "

# Row 28
$ws.Range("E28").Value = "  -0.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0927"
$ws.Range("E29").Value = "  -1.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.14"
$ws.Range("E30").Value = "  -1.56%  "

# Row 31
$ws.Range("E31").Value = "  -5.09%  "

# Row 32
$ws.Range("E32").Value = "  -1.30%  "

# Row 33
$ws.Range("E33").Value = "  -0.60%  "

# Row 34
$ws.Range("E34").Value = "  -4.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("E36").Value = "  -1.04%  "

# Row 37
$ws.Range("E37").Value = "  -1.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "153.38"
$ws.Range("E38").Value = "  +1.96%  "

# Row 39
$ws.Range("E39").Value = "  -4.97%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.48"
$ws.Range("E40").Value = "  -1.13%  "

# Row 41
$ws.Range("E41").Value = "  -1.78%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.50"
$ws.Range("E43").Value = "  -3.24%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.22"
$ws.Range("E44").Value = "  -5.09%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.62"
$ws.Range("E45").Value = "  -0.91%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0524"
$ws.Range("E46").Value = "  -3.37%  "

# Row 47 (was Mantle, becomes InjectiveProtocol)
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.76"
$ws.Range("E47").Value = "  -2.52%  "

# Row 48 (was InjectiveProtocol, becomes Mantle)
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.590"
$ws.Range("E48").Value = "  +0.05%  "

# Row 49
$ws.Range("E49").Value = "  -0.56%  "

# Row 50
$ws.Range("E50").Value = "  -1.24%  "

# Row 51
$ws.Range("E51").Value = "  +0.70%  "
